# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# Net content change (resolved from the OOXML diff):
#   - The "Periodo Mora" / "Valor Mora" pair on row 16 and row 18 are
#     swapped; row 17 is left untouched.
#       Before: E16=2106 F16=39480 | E17=2105 F17=39480 | E18=2104 F18=30268
#       After:  E16=2104 F16=30268 | E17=2105 F17=39480 | E18=2106 F18=39480

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E16").Value = "2104"
$ws.Range("F16").Value = 30268

$ws.Range("E18").Value = "2106"
$ws.Range("F18").Value = 39480
